$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted ahead of the existing data (row 569),
# pushing every subsequent row down by one (569->570, ..., 627->628).
$ws.Rows.Item(569).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A569").Value = 6
$ws.Range("B569").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C569").Value = "Metropolitana"
$ws.Range("D569").Value = 45194
$ws.Range("E569").Value = 13
$ws.Range("F569").Value = 100112032
$ws.Range("G569").Value = "Zapallo italiano"
$ws.Range("H569").Value = "Sin especificar"
$ws.Range("I569").Value = "Primera"
$ws.Range("J569").Value = 600
$ws.Range("K569").Value = 10000
$ws.Range("L569").Value = 12000
$ws.Range("M569").Value = 11167
$ws.Range("N569").Value = "`$/caja 50 unidades"
$ws.Range("O569").Value = "Región de Arica y Parinacota"
$ws.Range("P569").Value = 223
$ws.Range("Q569").Value = 50
$ws.Range("R569").Value = "Hortaliza"
